$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0 (ALC)
$ws.Range("H19").Value = 2450.25
$ws.Range("I19").Value = 1899
$ws.Range("K19").Value = 1899
$ws.Range("M19").Value = -1724

# Hunk 1 (ALC)
$ws.Range("H64").Value = 4984.375
$ws.Range("J64").Value = 4960
$ws.Range("L64").Value = 4960
$ws.Range("N64").Value = -5456

# Hunk 2 (ALC)
$ws.Range("H67").Value = 4984.375
$ws.Range("J67").Value = 4960
$ws.Range("L67").Value = 4960
$ws.Range("N67").Value = -6676

# Hunk 3 (ALC)
$ws.Range("H86").Value = 11196.333
$ws.Range("J86").Value = 11196.333
$ws.Range("L86").Value = 11196.333
$ws.Range("N86").Value = -13442.333

# Hunk 4 (ALC)
$ws.Range("H88").Value = 1169843.1
$ws.Range("J88").Value = 1364317
$ws.Range("L88").Value = 1364317
$ws.Range("N88").Value = -1365129

# Hunk 5 (ALC)
$ws.Range("H89").Value = 11196.333
$ws.Range("J89").Value = 11196.333
$ws.Range("L89").Value = 55981.665
$ws.Range("N89").Value = -67213.66500000001

# Hunk 6 (ALC)
$ws.Range("H91").Value = 1169843.1
$ws.Range("J91").Value = 1364317
$ws.Range("L91").Value = 1364317
$ws.Range("N91").Value = -1367125

# Hunk 7 (ALC)
$ws.Range("H92").Value = 855.3333
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Hunk 8 (ALC)
$ws.Range("H106").Value = 4747.625
$ws.Range("I106").Value = 1992.5
$ws.Range("J106").Value = 5666
$ws.Range("K106").Value = 1992.5
$ws.Range("L106").Value = 5666
$ws.Range("M106").Value = -1361.5
$ws.Range("N106").Value = -6928

# Hunk 9 (ALC)
$ws.Range("H113").Value = 3736.6667
$ws.Range("I113").Value = 3736.6667
$ws.Range("K113").Value = 3736.6667
$ws.Range("M113").Value = -482.6667000000002

# Hunk 10 (ALC)
$ws.Range("H125").Value = 2799.5
$ws.Range("I125").Value = 5000
$ws.Range("K125").Value = 45000
$ws.Range("M125").Value = -42540

$ws = $wb.Worksheets.Item("ARM")
# Hunk 11 (ARM)
$ws.Range("H44").Value = 39495
$ws.Range("J44").Value = 39495
$ws.Range("L44").Value = 39495
$ws.Range("N44").Value = -40471

# Hunk 12 (ARM)
$ws.Range("H45").Value = 3084.6667
$ws.Range("I45").Value = 1508
$ws.Range("J45").Value = 3400
$ws.Range("K45").Value = 1508
$ws.Range("L45").Value = 3400
$ws.Range("M45").Value = -1131
$ws.Range("N45").Value = -4154

# Hunk 13 (ARM)
$ws.Range("H110").Value = 4392
$ws.Range("I110").Value = 3990
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 3990
$ws.Range("L110").Value = 6000
$ws.Range("M110").Value = -1945
$ws.Range("N110").Value = -10090

# Hunk 14 (ARM)
$ws.Range("H122").Value = 2218.9
$ws.Range("I122").Value = 1597.25
$ws.Range("K122").Value = 4791.75
$ws.Range("M122").Value = -2341.75

# Hunk 15 (ARM)
$ws.Range("H132").Value = 11291.866
$ws.Range("I132").Value = 8670.522999999999
$ws.Range("J132").Value = 17408.334
$ws.Range("K132").Value = 26011.569
$ws.Range("L132").Value = 52225.00199999999
$ws.Range("M132").Value = -23481.569
$ws.Range("N132").Value = -57285.00199999999

$ws = $wb.Worksheets.Item("BSM")
# Hunk 16 (BSM)
$ws.Range("H20").Value = 5723.7
$ws.Range("I20").Value = 4807.875
$ws.Range("J20").Value = 9387
$ws.Range("K20").Value = 4807.875
$ws.Range("L20").Value = 9387
$ws.Range("M20").Value = -4560.875
$ws.Range("N20").Value = -9881

# Hunk 17 (BSM)
$ws.Range("H105").Value = 2131.25
$ws.Range("I105").Value = 872.8
$ws.Range("K105").Value = 872.8
$ws.Range("M105").Value = 874.2

# Hunk 18 (BSM)
$ws.Range("H107").Value = 1311.5625
$ws.Range("I107").Value = 1304.7142
$ws.Range("J107").Value = 1359.5
$ws.Range("K107").Value = 1304.7142
$ws.Range("L107").Value = 1359.5
$ws.Range("M107").Value = 615.2858000000001
$ws.Range("N107").Value = -5199.5

# Hunk 19 (BSM)
$ws.Range("H134").Value = 62633.445
$ws.Range("I134").Value = 3114.5
$ws.Range("K134").Value = 9343.5
$ws.Range("M134").Value = -6808.5

$ws = $wb.Worksheets.Item("CRP")
# Hunk 20 (CRP)
$ws.Range("H99").Value = 3784.25
$ws.Range("I99").Value = 3492.889
$ws.Range("J99").Value = 4658.3335
$ws.Range("K99").Value = 3492.889
$ws.Range("L99").Value = 4658.3335
$ws.Range("M99").Value = -1994.889
$ws.Range("N99").Value = -7654.3335

# Hunk 21 (CRP)
$ws.Range("H105").Value = 3701.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3701.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3701.5
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -7195.5

# Hunk 22 (CRP)
$ws.Range("H126").Value = 3784.25
$ws.Range("I126").Value = 3492.889
$ws.Range("J126").Value = 4658.3335
$ws.Range("K126").Value = 10478.667
$ws.Range("L126").Value = 13975.0005
$ws.Range("M126").Value = -8008.667000000001
$ws.Range("N126").Value = -18915.0005

# Hunk 23 (CRP)
$ws.Range("H132").Value = 2982.2
$ws.Range("I132").Value = 2980.3333
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 8940.999899999999
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -6410.999899999999
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("CUL")
# Hunk 24 (CUL)
$ws.Range("H68").Value = 2606
$ws.Range("J68").Value = 2460.5
$ws.Range("L68").Value = 7381.5
$ws.Range("N68").Value = -9003.5

# Hunk 25 (CUL)
$ws.Range("H71").Value = 2606
$ws.Range("J71").Value = 2460.5
$ws.Range("L71").Value = 22144.5
$ws.Range("N71").Value = -30256.5

# Hunk 26 (CUL)
$ws.Range("H125").Value = 12254.125
$ws.Range("J125").Value = 12254.125
$ws.Range("L125").Value = 36762.375
$ws.Range("N125").Value = -46602.375

# Hunk 27 (CUL)
$ws.Range("H131").Value = 7029.5835
$ws.Range("I131").Value = 8065.4287
$ws.Range("K131").Value = 24196.2861
$ws.Range("M131").Value = -19156.2861

$ws = $wb.Worksheets.Item("GSM")
# Hunk 28 (GSM)
$ws.Range("H102").Value = 6769.222
$ws.Range("I102").Value = 4653.8335
$ws.Range("J102").Value = 11000
$ws.Range("K102").Value = 4653.8335
$ws.Range("L102").Value = 11000
$ws.Range("M102").Value = -3031.8335
$ws.Range("N102").Value = -14244

# Hunk 29 (GSM)
$ws.Range("H113").Value = 3922.4285
$ws.Range("I113").Value = 3999.7778
$ws.Range("J113").Value = 3783.2
$ws.Range("K113").Value = 3999.7778
$ws.Range("L113").Value = 3783.2
$ws.Range("M113").Value = -1829.7778
$ws.Range("N113").Value = -8123.2

$ws = $wb.Worksheets.Item("LTW")
# Hunk 30 (LTW)
$ws.Range("H46").Value = 4145.364
$ws.Range("J46").Value = 4918.75
$ws.Range("L46").Value = 4918.75
$ws.Range("N46").Value = -5294.75

# Hunk 31 (LTW)
$ws.Range("H61").Value = 989
$ws.Range("I61").Value = 989
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 989
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -787
$ws.Range("N61").ClearContents()

# Hunk 32 (LTW)
$ws.Range("H113").Value = 989
$ws.Range("I113").Value = 989
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 989
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1181
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Hunk 33 (WVR)
$ws.Range("H24").Value = 30006.666
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 30006.666
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 30006.666
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -30466.666

# Hunk 34 (WVR)
$ws.Range("H113").Value = 805.1667
$ws.Range("I113").Value = 830
$ws.Range("J113").Value = 532
$ws.Range("K113").Value = 2490
$ws.Range("L113").Value = 1596
$ws.Range("M113").Value = -320
$ws.Range("N113").Value = -5936

# Hunk 35 (WVR)
$ws.Range("H132").Value = 10878.286
$ws.Range("I132").Value = 1465.5385
$ws.Range("K132").Value = 4396.6155
$ws.Range("M132").Value = -1866.6155
